$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-18 Thursday" "2024-01-19 Friday"

Replace-Text "218×5=1090" "766×6=4596"
Replace-Text "995×8=7960" "977×8=7816"
Replace-Text "354×2=708" "948×8=7584"
Replace-Text "762×3=2286" "330×4=1320"
Replace-Text "417×7=2919" "967×5=4835"

Replace-Text "669×6=4014" "710×7=4970"
Replace-Text "294×6=1764" "440×3=1320"
Replace-Text "375×9=3375" "130×6=780"
Replace-Text "691×4=2764" "520×6=3120"
Replace-Text "169×5=845" "360×4=1440"

Replace-Text "963×9=8667" "445×4=1780"
Replace-Text "829×4=3316" "543×5=2715"
Replace-Text "718×5=3590" "722×5=3610"
Replace-Text "381×3=1143" "618×8=4944"
Replace-Text "605×6=3630" "169×3=507"

Replace-Text "185×3=555" "821×8=6568"
Replace-Text "672×8=5376" "939×2=1878"
Replace-Text "438×9=3942" "871×3=2613"
Replace-Text "660×6=3960" "212×2=424"
Replace-Text "251×6=1506" "163×4=652"

Replace-Text "186×8=1488" "873×4=3492"
Replace-Text "379×5=1895" "908×5=4540"
Replace-Text "354×3=1062" "610×5=3050"
Replace-Text "912×5=4560" "805×9=7245"
Replace-Text "106×9=954" "798×5=3990"
